$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert three new rows at the top with documentation text
# (shared-string table order requires A2 to be written before A1)
$ws.Range("A2").Value = "referring to ICSolar.mo in github"
$ws.Range("A1").Value = "What are the receiver and tubing loss to ambient coefficients that we arrived at during RSME regression fitting?"
$ws.Range("A3").Value = "done for conf paper 2020-08-22"

# Apply a numeric format (0.000) to the final RC/summary results
$ws.Range("B41").NumberFormat = "0.000"
$ws.Range("B42").NumberFormat = "0.000"

# Update sheet view: clear prior topLeftCell/selection, select B30
$ws.Range("B30").Select()

# Page setup: explicit portrait orientation (adds a pageSetup element on save)
$ws.PageSetup.Orientation = 1
